$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 35294330
$ws.Range("I33").Value = 46154050
$ws.Range("K33").Value = 46154050
$ws.Range("M33").Value = -46153821
$ws.Range("H34").Value = 6324.222
$ws.Range("I34").Value = 5871.143
$ws.Range("J34").Value = 7910
$ws.Range("K34").Value = 5871.143
$ws.Range("L34").Value = 7910
$ws.Range("M34").Value = -5668.143
$ws.Range("N34").Value = -8316
$ws.Range("H36").Value = 6324.222
$ws.Range("I36").Value = 5871.143
$ws.Range("J36").Value = 7910
$ws.Range("K36").Value = 5871.143
$ws.Range("L36").Value = 7910
$ws.Range("M36").Value = -5156.143
$ws.Range("N36").Value = -9340
$ws.Range("H40").Value = 1248.3125
$ws.Range("J40").Value = 1441
$ws.Range("L40").Value = 1441
$ws.Range("N40").Value = -1791
$ws.Range("H43").Value = 10764.091
$ws.Range("I43").Value = 25737
$ws.Range("K43").Value = 25737
$ws.Range("M43").Value = -25668
$ws.Range("H96").Value = 312.83334
$ws.Range("I96").Value = 363.42856
$ws.Range("K96").Value = 1090.28568
$ws.Range("M96").Value = 282.71432
$ws.Range("H97").Value = 18366.77
$ws.Range("J97").Value = 42730.547
$ws.Range("L97").Value = 128191.641
$ws.Range("N97").Value = -129183.641
$ws.Range("H106").Value = 1480.5
$ws.Range("I106").Value = 875.625
$ws.Range("K106").Value = 875.625
$ws.Range("M106").Value = -244.625
$ws.Range("H112").Value = 1791.2142
$ws.Range("I112").Value = 933.6667
$ws.Range("K112").Value = 2801.0001
$ws.Range("M112").Value = -1693.0001
$ws.Range("H132").Value = 11178959
$ws.Range("I132").Value = 11941727
$ws.Range("K132").Value = 35825181
$ws.Range("M132").Value = -35822651
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 1625
$ws.Range("I35").Value = 833.3333
$ws.Range("J35").Value = 4000
$ws.Range("K35").Value = 833.3333
$ws.Range("L35").Value = 4000
$ws.Range("M35").Value = -427.3333
$ws.Range("N35").Value = -4812
$ws.Range("H56").Value = 9999
$ws.Range("I56").Value = 9999
$ws.Range("K56").Value = 9999
$ws.Range("M56").Value = -9257
$ws.Range("H61").Value = 4914.923
$ws.Range("I61").Value = 4535.909
$ws.Range("K61").Value = 4535.909
$ws.Range("M61").Value = -4323.909
$ws.Range("H122").Value = 3780
$ws.Range("I122").Value = 3780
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11340
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8890
$ws.Range("N122").Value = ""
$ws.Range("H136").Value = 4914.923
$ws.Range("I136").Value = 4535.909
$ws.Range("K136").Value = 13607.727
$ws.Range("M136").Value = -11057.727

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 500479
$ws.Range("J70").Value = 500479
$ws.Range("L70").Value = 500479
$ws.Range("N70").Value = -501065
$ws.Range("H73").Value = 500479
$ws.Range("J73").Value = 500479
$ws.Range("L73").Value = 500479
$ws.Range("N73").Value = -502507
$ws.Range("H86").Value = 6128.2285
$ws.Range("I86").Value = 6163.72
$ws.Range("J86").Value = 6039.5
$ws.Range("K86").Value = 6163.72
$ws.Range("L86").Value = 6039.5
$ws.Range("M86").Value = -5040.72
$ws.Range("N86").Value = -8285.5
$ws.Range("H89").Value = 6128.2285
$ws.Range("I89").Value = 6163.72
$ws.Range("J89").Value = 6039.5
$ws.Range("K89").Value = 30818.6
$ws.Range("L89").Value = 30197.5
$ws.Range("M89").Value = -25202.6
$ws.Range("N89").Value = -41429.5
$ws.Range("H99").Value = 3419.9333
$ws.Range("I99").Value = 2820.348
$ws.Range("K99").Value = 2820.348
$ws.Range("M99").Value = -1322.348
$ws.Range("H105").Value = 3429.1333
$ws.Range("I105").Value = 3535.2856
$ws.Range("K105").Value = 3535.2856
$ws.Range("M105").Value = -1788.2856
$ws.Range("H129").Value = 70780
$ws.Range("J129").Value = 70780
$ws.Range("L129").Value = 70780
$ws.Range("N129").Value = -80780
$ws.Range("H132").Value = 75698.57000000001
$ws.Range("J132").Value = 75698.57000000001
$ws.Range("L132").Value = 75698.57000000001
$ws.Range("N132").Value = -85818.57000000001
$ws.Range("H133").Value = 89995
$ws.Range("J133").Value = 89995
$ws.Range("L133").Value = 89995
$ws.Range("N133").Value = -100115
$ws.Range("H134").Value = 3684.2
$ws.Range("I134").Value = 3168.4
$ws.Range("K134").Value = 9505.200000000001
$ws.Range("M134").Value = -6970.200000000001
$ws.Range("H138").Value = 78860.664
$ws.Range("J138").Value = 78860.664
$ws.Range("L138").Value = 78860.664
$ws.Range("N138").Value = -89140.664
$ws.Range("H139").Value = 88000
$ws.Range("J139").Value = 88000
$ws.Range("L139").Value = 88000
$ws.Range("N139").Value = -98280

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5985477
$ws.Range("I31").Value = 7966589.5
$ws.Range("K31").Value = 7966589.5
$ws.Range("M31").Value = -7966294.5
$ws.Range("H34").Value = 5985477
$ws.Range("I34").Value = 7966589.5
$ws.Range("K34").Value = 7966589.5
$ws.Range("M34").Value = -7966387.5
$ws.Range("H42").Value = 4800
$ws.Range("I42").Value = 4800
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 4800
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -4207
$ws.Range("N42").Value = ""
$ws.Range("H80").Value = 13133.223
$ws.Range("J80").Value = 28799.334
$ws.Range("L80").Value = 28799.334
$ws.Range("N80").Value = -31045.334
$ws.Range("H83").Value = 13133.223
$ws.Range("J83").Value = 28799.334
$ws.Range("L83").Value = 86398.00199999999
$ws.Range("N83").Value = -97630.00199999999
$ws.Range("H107").Value = 1294.1666
$ws.Range("I107").Value = 1392.2222
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1392.2222
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 527.7778000000001
$ws.Range("N107").Value = -4840
$ws.Range("H122").Value = 1460.4286
$ws.Range("I122").Value = 1509.5454
$ws.Range("K122").Value = 4528.6362
$ws.Range("M122").Value = -2078.6362

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 71474.14
$ws.Range("I11").Value = 41.9
$ws.Range("J11").Value = 250054.75
$ws.Range("K11").Value = 125.7
$ws.Range("L11").Value = 750164.25
$ws.Range("M11").Value = 14.30000000000001
$ws.Range("N11").Value = -750444.25
$ws.Range("H26").Value = 1186.2174
$ws.Range("I26").Value = 133.07692
$ws.Range("J26").Value = 2555.3
$ws.Range("K26").Value = 399.23076
$ws.Range("L26").Value = 7665.900000000001
$ws.Range("M26").Value = -111.23076
$ws.Range("N26").Value = -8241.900000000001
$ws.Range("H113").Value = 534.94116
$ws.Range("I113").Value = 526.3333
$ws.Range("J113").Value = 599.5
$ws.Range("K113").Value = 1578.9999
$ws.Range("L113").Value = 1798.5
$ws.Range("M113").Value = 591.0001
$ws.Range("N113").Value = -6138.5
$ws.Range("H131").Value = 5145.615
$ws.Range("I131").Value = 2644
$ws.Range("K131").Value = 7932
$ws.Range("M131").Value = -2892
$ws.Range("H139").Value = 4215.4736
$ws.Range("I139").Value = 3630.875
$ws.Range("K139").Value = 10892.625
$ws.Range("M139").Value = -5752.625
$ws.Range("H140").Value = 8076.6
$ws.Range("I140").Value = 8707.764999999999
$ws.Range("K140").Value = 26123.295
$ws.Range("M140").Value = -20943.295

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 21481.883
$ws.Range("I126").Value = 28266
$ws.Range("K126").Value = 84798
$ws.Range("M126").Value = -82328
$ws.Range("H132").Value = 297291.53
$ws.Range("J132").Value = 4318.6665
$ws.Range("L132").Value = 12955.9995
$ws.Range("N132").Value = -18015.9995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4257.3125
$ws.Range("I22").Value = 1141.9
$ws.Range("J22").Value = 9449.666999999999
$ws.Range("K22").Value = 1141.9
$ws.Range("L22").Value = 9449.666999999999
$ws.Range("M22").Value = -846.9000000000001
$ws.Range("N22").Value = -10039.667
$ws.Range("H27").Value = 4257.3125
$ws.Range("I27").Value = 1141.9
$ws.Range("J27").Value = 9449.666999999999
$ws.Range("K27").Value = 1141.9
$ws.Range("L27").Value = 9449.666999999999
$ws.Range("M27").Value = -1034.9
$ws.Range("N27").Value = -9663.666999999999
$ws.Range("H46").Value = 3943.7144
$ws.Range("I46").Value = 1921.7
$ws.Range("J46").Value = 8998.75
$ws.Range("K46").Value = 1921.7
$ws.Range("L46").Value = 8998.75
$ws.Range("M46").Value = -1733.7
$ws.Range("N46").Value = -9374.75
$ws.Range("H132").Value = 6760
$ws.Range("I132").Value = 5620
$ws.Range("K132").Value = 16860
$ws.Range("M132").Value = -14330
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = ""

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
